# budget_project.xlsx - "modifications on the budget"
#
# Applies the three changes captured in the source-control diff:
#   1. Workbook window position reset to the top-left corner of the screen
#      (xWindow/yWindow -> 0/0 on the saved <workbookView>).
#   2. The "15 k $" label in the financial annex is reworded/re-ordered to
#      "$ 15 k " (leading currency sign, trailing space preserved).
#   3. The sheet's scroll position/selection moves up one row: the view
#      scrolls so row 4 is back at the top of the window and the active
#      cell changes from B10 to B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Reset the workbook window to the top-left of the screen ------------
# (xWindow="0" yWindow="0" in bookViews/workbookView)
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 0
$wb.Windows.Item(1).Left = 0
$wb.Windows.Item(1).Top = 0

# --- 2. Reword the "15 k $" budget cell to "$ 15 k " ------------------------
$ws.Range("B9").Value = "$ 15 k "

# --- 3. Scroll the sheet view up and move the selection to B9 --------------
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B9").Select()
